$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44326, 39, 211, 111.6325332120013),
    @(44327, 12, 209, 110.5744049351103),
    @(44328, 13, 211, 111.6325332120013),
    @(44329, 50, 223, 117.9813028733473)
)

$startRow = 252
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]

    # Match the style of column A used by the preceding rows (date format + border)
    $ws.Cells.Item(251, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
